$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural changes -------------------------------------------------
# Insert two new rows at the top (metadata: TsVersion / sourcelanguage /
# language header row, plus the "2.1" value row). Everything below shifts
# down by 2 rows.
$ws.Rows("1:2").Insert()

# Insert a new column B for "ID" (between Context and Source).
# Old Source/Translation/... columns shift one to the right.
$ws.Columns("B:B").Insert()

# Insert four new columns E:H for TranslationType / comment / extracomment /
# translatorcomment (between Translation and Location). The old Location
# columns shift four to the right.
$ws.Columns("E:H").Insert()

# --- New metadata rows ---------------------------------------------------
$ws.Range("A1").Value = "TsVersion"
$ws.Range("B1").Value = "sourcelanguage"
$ws.Range("C1").Value = "language"

# "2.1" is the TS format version, stored as text (not a number) in the
# source file, so force a text number format before writing it.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2.1"

# --- Header row (was row 1, now row 3) -----------------------------------
$ws.Range("A3").Value = "Context"
$ws.Range("B3").Value = "ID"
$ws.Range("C3").Value = "Source"
$ws.Range("D3").Value = "Translation"
$ws.Range("E3").Value = "TranslationType"
$ws.Range("F3").Value = "comment"
$ws.Range("G3").Value = "extracomment"
$ws.Range("H3").Value = "translatorcomment"
$ws.Range("I3").Value = "Location"

# --- Data rows (were rows 2-7, now rows 4-9) ------------------------------
$ws.Range("A4").Value = "ThemeWidget"
$ws.Range("C4").Value = "Series"
$ws.Range("I4").Value = "../themewidget.cpp - 289"
$ws.Range("J4").Value = "../themewidget.cpp - 290"
$ws.Range("K4").Value = "../themewidget.cpp - 291"

$ws.Range("A5").Value = "Test"
$ws.Range("C5").Value = "SourceValue"
$ws.Range("I5").Value = "../themewidget.cpp - 89"
$ws.Range("J5").Value = "../themewidget.cpp - 90"
$ws.Range("K5").Value = "../themewidget.cpp - 91"

$ws.Range("A6").Value = "ThemeWidgetForm"
$ws.Range("C6").Value = "Theme:"
$ws.Range("I6").Value = "../themewidget.ui - 19"

$ws.Range("A7").Value = "ThemeWidgetForm"
$ws.Range("C7").Value = "Animation:"
$ws.Range("I7").Value = "../themewidget.ui - 29"

$ws.Range("A8").Value = "ThemeWidgetForm"
$ws.Range("C8").Value = "Legend:"
$ws.Range("I8").Value = "../themewidget.ui - 39"

$ws.Range("A9").Value = "ThemeWidgetForm"
$ws.Range("C9").Value = "Anti-aliasing"
$ws.Range("I9").Value = "../themewidget.ui - 49"
